$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole data range to Text format first so that numeric-looking
# strings (prices, percentages, the hour marker, etc.) are preserved exactly
# as text instead of being auto-converted to numbers/percentages by Excel.
$dataRange = $ws.Range("B2:G51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "290.60"
$ws.Range("E2").Value = "-4.13%"
$ws.Range("G2").Value = "4"
$ws.Range("D3").Value = "30.79"
$ws.Range("E3").Value = "-4.29%"
$ws.Range("G3").Value = "4"
$ws.Range("D4").Value = "4.958"
$ws.Range("E4").Value = "0.44%"
$ws.Range("G4").Value = "4"
$ws.Range("E5").Value = "-8.62%"
$ws.Range("G5").Value = "4"
$ws.Range("D6").Value = "1.783"
$ws.Range("E6").Value = "-13.46%"
$ws.Range("G6").Value = "4"
$ws.Range("D7").Value = "7.660"
$ws.Range("E7").Value = "-2.19%"
$ws.Range("G7").Value = "4"
$ws.Range("D8").Value = "3.735"
$ws.Range("E8").Value = "-3.03%"
$ws.Range("G8").Value = "4"
$ws.Range("D9").Value = "0.8959"
$ws.Range("E9").Value = "-3.24%"
$ws.Range("G9").Value = "4"
$ws.Range("D10").Value = "0.1643"
$ws.Range("E10").Value = "-6.60%"
$ws.Range("G10").Value = "4"
$ws.Range("D11").Value = "0.07692"
$ws.Range("E11").Value = "-1.89%"
$ws.Range("G11").Value = "4"
$ws.Range("D12").Value = "0.08007"
$ws.Range("E12").Value = "-7.05%"
$ws.Range("G12").Value = "4"
$ws.Range("D13").Value = "0.03026"
$ws.Range("E13").Value = "-4.00%"
$ws.Range("G13").Value = "4"
$ws.Range("E14").Value = "-0.26%"
$ws.Range("G14").Value = "4"
$ws.Range("D15").Value = "0.001503"
$ws.Range("E15").Value = "-0.98%"
$ws.Range("G15").Value = "4"
$ws.Range("D16").Value = "0.005733"
$ws.Range("E16").Value = "-0.27%"
$ws.Range("G16").Value = "4"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.469"
$ws.Range("E17").Value = "0.13%"
$ws.Range("G17").Value = "4"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.080"
$ws.Range("E18").Value = "-1.82%"
$ws.Range("G18").Value = "4"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3278"
$ws.Range("E19").Value = "0.06%"
$ws.Range("G19").Value = "4"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1330"
$ws.Range("E20").Value = "0.80%"
$ws.Range("G20").Value = "4"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "4.046"
$ws.Range("E21").Value = "-5.41%"
$ws.Range("G21").Value = "4"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.1999"
$ws.Range("E22").Value = "0.35%"
$ws.Range("G22").Value = "4"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "0.04515"
$ws.Range("E23").Value = "-1.32%"
$ws.Range("G23").Value = "4"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "-0.80%"
$ws.Range("G24").Value = "4"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "0.004004"
$ws.Range("E25").Value = "-9.96%"
$ws.Range("G25").Value = "4"
$ws.Range("B26").Value = "NitroEx"
$ws.Range("C26").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D26").Value = "0.0001250"
$ws.Range("E26").Value = "-0.13%"
$ws.Range("G26").Value = "4"
$ws.Range("B27").Value = "Spectre.aiUtilityToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("D27").Value = "--"
$ws.Range("E27").Value = "--%"
$ws.Range("G27").Value = "4"
$ws.Range("B28").Value = "LegolasExchange"
$ws.Range("C28").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("G28").Value = "4"
$ws.Range("B29").Value = "BitZToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("G29").Value = "4"
$ws.Range("B30").Value = "Birake"
$ws.Range("C30").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("G30").Value = "4"
$ws.Range("B31").Value = "NashExchange"
$ws.Range("C31").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("G31").Value = "4"
$ws.Range("B32").Value = "AAXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("G32").Value = "4"
$ws.Range("B33").Value = "CenX"
$ws.Range("C33").Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("G33").Value = "4"
$ws.Range("B34").Value = "BNIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("G34").Value = "4"
$ws.Range("B35").Value = "UpBots"
$ws.Range("C35").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("G35").Value = "4"
$ws.Range("G36").Value = "4"
$ws.Range("G37").Value = "4"
$ws.Range("G38").Value = "4"
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").Value = "-7.04%"
$ws.Range("G39").Value = "4"
$ws.Range("D40").Value = "0.04363"
$ws.Range("E40").Value = "-9.34%"
$ws.Range("G40").Value = "4"
$ws.Range("D41").Value = "0.007384"
$ws.Range("E41").Value = "-1.57%"
$ws.Range("G41").Value = "4"
$ws.Range("D42").Value = "0.1307"
$ws.Range("E42").Value = "-4.07%"
$ws.Range("G42").Value = "4"
$ws.Range("D43").Value = "0.002059"
$ws.Range("E43").Value = "-12.83%"
$ws.Range("G43").Value = "4"
$ws.Range("D44").Value = "0.009232"
$ws.Range("E44").Value = "-12.52%"
$ws.Range("G44").Value = "4"
$ws.Range("D45").Value = "0.00005975"
$ws.Range("E45").Value = "-5.64%"
$ws.Range("G45").Value = "4"
$ws.Range("E46").Value = "-0.12%"
$ws.Range("G46").Value = "4"
$ws.Range("G47").Value = "4"
$ws.Range("D48").Value = "0.002999"
$ws.Range("E48").Value = "-3.30%"
$ws.Range("G48").Value = "4"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").Value = "-0.12%"
$ws.Range("G49").Value = "4"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").Value = "-0.12%"
$ws.Range("G50").Value = "4"
$ws.Range("G51").Value = "4"

# Restore the default "Normal" style so no stray number-format styling is left behind
$dataRange.Style = "Normal"
